# live_trading_results.xlsx - apply trade #59 close + trade #116 open updates
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh aggregate metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.12   # Current Capital
$summary.Range("B4").Value = 0.23      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 87        # Total Trades
$summary.Range("B8").Value = 35        # Losing Trades
$summary.Range("B9").Value = 49.43     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": refresh the "momentum" strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.19
$status.Range("D11").Value = 17
$status.Range("E11").Value = -0.8100000000000001
$status.Range("F11").Value = -0.8100000000000001
$status.Range("G11").Value = 17.65

# ---------------------------------------------------------------------------
# Sheet "All Trades": trade #87 (momentum) closes, trade #116 (MarketMaking) opens
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #87 (row 88) transitions from OPEN to CLOSED
$allTrades.Range("G88").Value = 0.9399999999999999
$allTrades.Range("H88").Value = "CLOSED"
$allTrades.Range("I88").Value = -2.0833
$allTrades.Range("J88").Value = -0.02
$allTrades.Range("K88").Value = 99.19
$allTrades.Range("L88").Value = "early_exit"
$allTrades.Range("M88").Value = 0.14

# New trade #116 (row 117) - MarketMaking, newly opened
$allTrades.Range("A117").Value = 116
$allTrades.Range("B117:C117").NumberFormat = "@"
$allTrades.Range("B117").Value = "2026-02-18"
$allTrades.Range("C117").Value = "00:22:59"
$allTrades.Range("D117").Value = "MarketMaking"
$allTrades.Range("E117").Value = "DOWN"
$allTrades.Range("F117").Value = 0.96
$allTrades.Range("H117").Value = "OPEN"
$allTrades.Range("I117").Value = 0
$allTrades.Range("J117").Value = 0
$allTrades.Range("K117").Value = 99.410254715139
$allTrades.Range("M117").Value = 0
$allTrades.Range("N117").Value = 0
$allTrades.Range("O117").Value = 0
$allTrades.Range("P117").Value = 0.6
$allTrades.Range("Q117").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# Sheet "momentum": trade #87 (row 18) transitions from OPEN to CLOSED
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("G18").Value = 0.9399999999999999
$momentum.Range("H18").Value = "CLOSED"
$momentum.Range("I18").Value = -2.0833
$momentum.Range("J18").Value = -0.02
$momentum.Range("K18").Value = 99.19
$momentum.Range("P18").Value = "early_exit"
$momentum.Range("Q18").Value = 0.14

# ---------------------------------------------------------------------------
# Sheet "MarketMaking": new trade #116 (row 37) - newly opened
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A37").Value = 116
$marketMaking.Range("B37:C37").NumberFormat = "@"
$marketMaking.Range("B37").Value = "2026-02-18"
$marketMaking.Range("C37").Value = "00:22:59"
$marketMaking.Range("D37").Value = "MarketMaking"
$marketMaking.Range("E37").Value = "DOWN"
$marketMaking.Range("F37").Value = 0.96
$marketMaking.Range("H37").Value = "OPEN"
$marketMaking.Range("I37").Value = 0
$marketMaking.Range("J37").Value = 0
$marketMaking.Range("K37").Value = 99.410254715139
$marketMaking.Range("L37").Value = 0
$marketMaking.Range("M37").Value = 0
$marketMaking.Range("N37").Value = 0.6
$marketMaking.Range("O37").Value = "Normal spread capture: 198 bps"
$marketMaking.Range("Q37").Value = 0
